$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "51.209.28"
$ws.Range("E2").Value = "  -0.74%  "
$ws.Range("D3").Value = "3.063.55"
$ws.Range("E3").Value = "  +1.29%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'388.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.27%  "
$ws.Range("D6").Value = "'102.02"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.15%  "
$ws.Range("E7").Value = "  -1.70%  "
$ws.Range("E8").Value = "  +0.01%  "
$ws.Range("D9").Value = "'0.578"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.64%  "
$ws.Range("D10").Value = "'36.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.23%  "
$ws.Range("E11").Value = "  +0.27%  "
$ws.Range("E12").Value = "  -1.19%  "
$ws.Range("D13").Value = "3.545.74"
$ws.Range("E13").Value = "  +1.23%  "
$ws.Range("D14").Value = "'18.26"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").Value = "'7.68"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.59%  "
$ws.Range("D16").Value = "3.059.48"
$ws.Range("E16").Value = "  +0.68%  "
$ws.Range("D17").Value = "'0.994"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.09%  "
$ws.Range("D18").Value = "'10.68"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.80%  "
$ws.Range("D19").Value = "51.220.84"
$ws.Range("E19").Value = "  -0.72%  "
$ws.Range("D20").Value = "'3.19"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.10%  "
$ws.Range("D21").Value = "'12.25"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.30%  "
$ws.Range("D22").Value = "0.0₃0954"
$ws.Range("D23").Value = "'69.68"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.33%  "
$ws.Range("D24").Value = "'264.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -1.07%  "
$ws.Range("E25").Value = "  -1.15%  "
$ws.Range("E26").Value = "  -5.89%  "
$ws.Range("D27").Value = "'26.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.10%  "
$ws.Range("D28").Value = "'7.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -3.24%  "
$ws.Range("E30").Value = "  -5.51%  "
$ws.Range("E31").Value = "  -3.32%  "
$ws.Range("E32").Value = "  +1.56%  "
$ws.Range("D33").Value = "'35.67"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.93%  "
$ws.Range("D34").Value = "'0.0472"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +4.86%  "
$ws.Range("E35").Value = "  +1.41%  "
$ws.Range("D36").Value = "'50.06"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.02%  "
$ws.Range("D37").Value = "'1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.06%  "
$ws.Range("D38").Value = "'3.34"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.97%  "
$ws.Range("D39").Value = "'0.294"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.67%  "
$ws.Range("D40").Value = "'130.37"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.54%  "
$ws.Range("B41").Value = "Celestia"
$ws.Range("C41").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D41").Value = "'16.50"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.09%  "
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").Value = "'1.82"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.42%  "
$ws.Range("E43").Value = "  -0.77%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'2.47"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.07%  "
$ws.Range("B45").Value = "NEARProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D45").Value = "'3.73"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").Value = "'21.67"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("E47").Value = "  +3.47%  "
$ws.Range("E48").Value = "  -0.20%  "
$ws.Range("D49").Value = "2.063.07"
$ws.Range("E49").Value = "  +2.20%  "
$ws.Range("D50").Value = "'0.0325"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.80%  "
$ws.Range("D51").Value = "'0.905"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +14.52%  "
